# fix latency units in report sheets
# X2 header: "Utility" -> "Utility (Percent)"
# For each data row (3-23): append " msec"/" usec" to the Read Latency
# triple (L,M,N) and the Write Latency triple (O,P,Q) based on the
# original magnitude-derived unit recorded per column below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('X2').Value = 'Utility (Percent)'

# row => @{ base values for L,M,N,O,P,Q; read/write unit }
$rows = @(
    @{ Row=3; L='14'; M='2482'; N='1259.97'; O='1378'; P='2966.9k'; Q='45068.10'; ReadUnit='msec'; WriteUnit='usec' }
    @{ Row=4; L='22'; M='2985'; N='1257.68'; O='1398'; P='3559.6k'; Q='50624.71'; ReadUnit='msec'; WriteUnit='usec' }
    @{ Row=5; L='575'; M='10279k'; N='1209703.32'; O='1415'; P='11024k'; Q='349783.45'; ReadUnit='usec'; WriteUnit='usec' }
    @{ Row=6; L='10'; M='2500'; N='1276.03'; O='1368'; P='565907'; Q='7619.81'; ReadUnit='msec'; WriteUnit='usec' }
    @{ Row=7; L='10'; M='3743'; N='1212.89'; O='1421'; P='4460.0k'; Q='156258.79'; ReadUnit='msec'; WriteUnit='usec' }
    @{ Row=8; L='2'; M='2507'; N='1266.32'; O='1267'; P='2499.7k'; Q='30571.82'; ReadUnit='msec'; WriteUnit='usec' }
    @{ Row=9; L='10'; M='2719'; N='1275.16'; O='1370'; P='2111.7k'; Q='9532.50'; ReadUnit='msec'; WriteUnit='usec' }
    @{ Row=10; L='559'; M='5710.2k'; N='1236817.47'; O='1326'; P='7710.7k'; Q='124094.45'; ReadUnit='usec'; WriteUnit='usec' }
    @{ Row=11; L='4'; M='3365'; N='1277.56'; O='1326'; P='2150.2k'; Q='3971.32'; ReadUnit='msec'; WriteUnit='usec' }
    @{ Row=12; L='13'; M='3256'; N='1270.45'; O='2'; P='3003'; Q='20.52'; ReadUnit='msec'; WriteUnit='msec' }
    @{ Row=13; L='1271'; M='3214.9k'; N='1175223.88'; O='1445'; P='3773.0k'; Q='243746.25'; ReadUnit='usec'; WriteUnit='usec' }
    @{ Row=14; L='10'; M='2943'; N='1270.74'; O='2'; P='2248'; Q='19.71'; ReadUnit='msec'; WriteUnit='msec' }
    @{ Row=15; L='298'; M='5972.5k'; N='775063.11'; O='2'; P='9009'; Q='1524.04'; ReadUnit='usec'; WriteUnit='msec' }
    @{ Row=16; L='265'; M='5431.8k'; N='684835.57'; O='10'; P='6647'; Q='1723.49'; ReadUnit='usec'; WriteUnit='msec' }
    @{ Row=17; L='1079'; M='2487.5k'; N='1267433.54'; O='1264'; P='3352.6k'; Q='27820.85'; ReadUnit='usec'; WriteUnit='usec' }
    @{ Row=18; L='12'; M='2479'; N='1276.81'; O='1438'; P='1010.8k'; Q='5828.10'; ReadUnit='msec'; WriteUnit='usec' }
    @{ Row=19; L='5'; M='2745'; N='1272.97'; O='1557'; P='2135.2k'; Q='15084.99'; ReadUnit='msec'; WriteUnit='usec' }
    @{ Row=20; L='23'; M='2565'; N='1271.29'; O='2'; P='1696'; Q='18.75'; ReadUnit='msec'; WriteUnit='msec' }
    @{ Row=21; L='283'; M='5212.2k'; N='687479.44'; O='2'; P='6448'; Q='1742.68'; ReadUnit='usec'; WriteUnit='msec' }
    @{ Row=22; L='10'; M='2516'; N='1275.74'; O='1290'; P='991327'; Q='8301.50'; ReadUnit='msec'; WriteUnit='usec' }
    @{ Row=23; L='12'; M='2716'; N='1273.53'; O='1409'; P='2127.7k'; Q='13573.38'; ReadUnit='msec'; WriteUnit='usec' }
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Range("L$r").Value = $row.L + " " + $row.ReadUnit
    $ws.Range("M$r").Value = $row.M + " " + $row.ReadUnit
    $ws.Range("N$r").Value = $row.N + " " + $row.ReadUnit
    $ws.Range("O$r").Value = $row.O + " " + $row.WriteUnit
    $ws.Range("P$r").Value = $row.P + " " + $row.WriteUnit
    $ws.Range("Q$r").Value = $row.Q + " " + $row.WriteUnit
}

